$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (G1) onto the
# new H1 header cell, then set its text - this reuses the same style
# record instead of creating a near-duplicate one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
